$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The course table(s) got two blank rows inserted above them (rows 1-2),
# pushing all existing rows down by two (old row 1 -> row 3, ... old row
# 32 -> row 34). Shift everything down by inserting 2 new rows at the top.
$ws.Rows("1:2").Insert() | Out-Null

# Update the view state to match: the active selection moves to B20 and
# the window is scrolled down so row 4 is the top visible row.
$ws.Range("B20").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4

